# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the "handback" report
# generation pass: the status text moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language sheets get their
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (with a matching hyperlink for the new
# "Latest Target File" entries), and a few columns are widened so the
# newly-populated long file names are readable.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# -----------------------------------------------------------------------
# 1. Status text update (shared by both rows on the Overview sheet).
# -----------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# -----------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (I), "Latest Handback
#    File" (J) and "Latest Handback DateTime" (K) for both rows, and add
#    a hyperlink on the newly-populated "Latest Target File" cells
#    (pointing at the same source doc as column A's hyperlink).
# -----------------------------------------------------------------------
$ws2.Range("I2").Value = "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md"
$ws2.Range("J2").Value = "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.941f60b1ef4279cc6afb6cfeba2f4bc4a85ebe88.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-18 05:01:41"

$ws2.Range("I3").Value = "addcc77c-c26d-4265-8dcd-e95694c53179.md"
$ws2.Range("J3").Value = "addcc77c-c26d-4265-8dcd-e95694c53179.93f7c4f55f75202110437b2ce4a7992cfb42dc4a.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-18 05:01:41"

$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09e8e57d0799b5d49107becc91643722c2219db2/e2e/7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md", "", "", "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09e8e57d0799b5d49107becc91643722c2219db2/e2e/addcc77c-c26d-4265-8dcd-e95694c53179.md", "", "", "addcc77c-c26d-4265-8dcd-e95694c53179.md")

# -----------------------------------------------------------------------
# 3. de-de sheet: same treatment as zh-cn, but with the de-de xliff
#    names / handback timestamp.
# -----------------------------------------------------------------------
$ws3.Range("I2").Value = "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md"
$ws3.Range("J2").Value = "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.941f60b1ef4279cc6afb6cfeba2f4bc4a85ebe88.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-18 05:01:48"

$ws3.Range("I3").Value = "addcc77c-c26d-4265-8dcd-e95694c53179.md"
$ws3.Range("J3").Value = "addcc77c-c26d-4265-8dcd-e95694c53179.93f7c4f55f75202110437b2ce4a7992cfb42dc4a.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-18 05:01:48"

$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09e8e57d0799b5d49107becc91643722c2219db2/e2e/7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md", "", "", "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09e8e57d0799b5d49107becc91643722c2219db2/e2e/addcc77c-c26d-4265-8dcd-e95694c53179.md", "", "", "addcc77c-c26d-4265-8dcd-e95694c53179.md")

# -----------------------------------------------------------------------
# 4. Widen columns that now hold the long file names / are displayed
#    wider in the handback report.
#    (ColumnWidth is quantized to whole "characters" by this engine, so
#    the inputs below are chosen to land as close as possible to the
#    target stored widths.)
# -----------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668   # E
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668   # F

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668   # C
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664   # I
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664  # J

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668   # C
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664   # I
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664  # J
